$wb = $excel.ActiveWorkbook

$wsExhibition = $wb.Worksheets.Item("展览")
$wsExhibition.Range("F3").Value = 7574
$wsExhibition.Range("F5").Value = 19
$wsExhibition.Range("F7").Value = 4197
$wsExhibition.Range("F9").Value = 583
$wsExhibition.Range("F11").Value = 670

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value = 7574
$wsAll.Range("F7").Value = 19
$wsAll.Range("F9").Value = 4197
$wsAll.Range("F11").Value = 583
$wsAll.Range("F13").Value = 670
